$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 103,
# pushing every subsequent record (old rows 103-123) down by one row
# (new rows 104-124). Insert a row at 103 to replicate that shift while
# keeping formatting (e.g. the date style on column D) consistent with
# the rest of the table.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new record.
$ws.Range("A103").Value = 5
$ws.Range("B103").Value = "Macroferia Regional de Talca"
$ws.Range("C103").Value = "Maule"
$ws.Range("D103").Value = 45204
$ws.Range("E103").Value = 7
$ws.Range("F103").Value = 100112026
$ws.Range("G103").Value = "Haba"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 10000
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = 10000
$ws.Range("N103").Value = "$/saco 25 kilos"
$ws.Range("O103").Value = "Región de O'Higgins"
$ws.Range("P103").Value = 400
$ws.Range("Q103").Value = 25
$ws.Range("R103").Value = "Hortaliza"
